$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "S100a8"
$ws.Cells.Item(2, 3).Value = "Tlr4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.213427
$ws.Cells.Item(2, 8).Value = 0.640281
$ws.Cells.Item(2, 9).Value = 0.05929937785206704
$ws.Cells.Item(2, 10).Value = 0.05929937785206704
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 7.407905
$ws.Cells.Item(2, 14).Value = 22.223715
$ws.Cells.Item(2, 15).Value = 0.1577242380174723
$ws.Cells.Item(2, 16).Value = 0.1577242380174723
$ws.Cells.Item(2, 17).Value = 1.581046940435
$ws.Cells.Item(2, 18).Value = 14.229422463915
$ws.Cells.Item(2, 19).Value = 0.009352949186627447
$ws.Cells.Item(2, 20).Value = 0.009352949186627449

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "S100a8"
$ws.Cells.Item(3, 3).Value = "Tlr4"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.213427
$ws.Cells.Item(3, 8).Value = 0.640281
$ws.Cells.Item(3, 9).Value = 0.05929937785206704
$ws.Cells.Item(3, 10).Value = 0.05929937785206704
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 16.48752133333333
$ws.Cells.Item(3, 14).Value = 49.462564
$ws.Cells.Item(3, 15).Value = 0.3510414535684271
$ws.Cells.Item(3, 16).Value = 0.3510414535684271
$ws.Cells.Item(3, 17).Value = 3.518882215609334
$ws.Cells.Item(3, 18).Value = 31.669939940484
$ws.Cells.Item(3, 19).Value = 0.020816539796893
$ws.Cells.Item(3, 20).Value = 0.02081653979689301

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "S100a8"
$ws.Cells.Item(4, 3).Value = "Tlr4"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.213427
$ws.Cells.Item(4, 8).Value = 0.640281
$ws.Cells.Item(4, 9).Value = 0.05929937785206704
$ws.Cells.Item(4, 10).Value = 0.05929937785206704
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 5.464566666666666
$ws.Cells.Item(4, 14).Value = 16.3937
$ws.Cells.Item(4, 15).Value = 0.1163479571613943
$ws.Cells.Item(4, 16).Value = 0.1163479571613943
$ws.Cells.Item(4, 17).Value = 1.166286069966667
$ws.Cells.Item(4, 18).Value = 10.4965746297
$ws.Cells.Item(4, 19).Value = 0.006899361474029629
$ws.Cells.Item(4, 20).Value = 0.006899361474029629

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "S100a8"
$ws.Cells.Item(5, 3).Value = "Tlr4"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.213427
$ws.Cells.Item(5, 8).Value = 0.640281
$ws.Cells.Item(5, 9).Value = 0.05929937785206704
$ws.Cells.Item(5, 10).Value = 0.05929937785206704
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 17.60745533333333
$ws.Cells.Item(5, 14).Value = 52.822366
$ws.Cells.Item(5, 15).Value = 0.3748863512527063
$ws.Cells.Item(5, 16).Value = 0.3748863512527063
$ws.Cells.Item(5, 17).Value = 3.757906369427334
$ws.Cells.Item(5, 18).Value = 33.821157324846
$ws.Cells.Item(5, 19).Value = 0.02223052739451695
$ws.Cells.Item(5, 20).Value = 0.02223052739451696

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "S100a8"
$ws.Cells.Item(6, 3).Value = "Tlr4"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 1.115191333333333
$ws.Cells.Item(6, 8).Value = 3.345574
$ws.Cells.Item(6, 9).Value = 0.3098490455878768
$ws.Cells.Item(6, 10).Value = 0.3098490455878768
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 7.407905
$ws.Cells.Item(6, 14).Value = 22.223715
$ws.Cells.Item(6, 15).Value = 0.1577242380174723
$ws.Cells.Item(6, 16).Value = 0.1577242380174723
$ws.Cells.Item(6, 17).Value = 8.261231454156666
$ws.Cells.Item(6, 18).Value = 74.35108308740999
$ws.Cells.Item(6, 19).Value = 0.0488707046157889
$ws.Cells.Item(6, 20).Value = 0.04887070461578891

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "S100a8"
$ws.Cells.Item(7, 3).Value = "Tlr4"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 1.115191333333333
$ws.Cells.Item(7, 8).Value = 3.345574
$ws.Cells.Item(7, 9).Value = 0.3098490455878768
$ws.Cells.Item(7, 10).Value = 0.3098490455878768
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 16.48752133333333
$ws.Cells.Item(7, 14).Value = 49.462564
$ws.Cells.Item(7, 15).Value = 0.3510414535684271
$ws.Cells.Item(7, 16).Value = 0.3510414535684271
$ws.Cells.Item(7, 17).Value = 18.38674089908178
$ws.Cells.Item(7, 18).Value = 165.480668091736
$ws.Cells.Item(7, 19).Value = 0.1087698593499581
$ws.Cells.Item(7, 20).Value = 0.1087698593499581

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "S100a8"
$ws.Cells.Item(8, 3).Value = "Tlr4"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 1.115191333333333
$ws.Cells.Item(8, 8).Value = 3.345574
$ws.Cells.Item(8, 9).Value = 0.3098490455878768
$ws.Cells.Item(8, 10).Value = 0.3098490455878768
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 5.464566666666666
$ws.Cells.Item(8, 14).Value = 16.3937
$ws.Cells.Item(8, 15).Value = 0.1163479571613943
$ws.Cells.Item(8, 16).Value = 0.1163479571613943
$ws.Cells.Item(8, 17).Value = 6.094037387088889
$ws.Cells.Item(8, 18).Value = 54.84633648379999
$ws.Cells.Item(8, 19).Value = 0.03605030348255719
$ws.Cells.Item(8, 20).Value = 0.03605030348255719

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "S100a8"
$ws.Cells.Item(9, 3).Value = "Tlr4"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 1.115191333333333
$ws.Cells.Item(9, 8).Value = 3.345574
$ws.Cells.Item(9, 9).Value = 0.3098490455878768
$ws.Cells.Item(9, 10).Value = 0.3098490455878768
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 17.60745533333333
$ws.Cells.Item(9, 14).Value = 52.822366
$ws.Cells.Item(9, 15).Value = 0.3748863512527063
$ws.Cells.Item(9, 16).Value = 0.3748863512527063
$ws.Cells.Item(9, 17).Value = 19.63568158978712
$ws.Cells.Item(9, 18).Value = 176.721134308084
$ws.Cells.Item(9, 19).Value = 0.1161581781395726
$ws.Cells.Item(9, 20).Value = 0.1161581781395726

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "S100a8"
$ws.Cells.Item(10, 3).Value = "Tlr4"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.3440983333333333
$ws.Cells.Item(10, 8).Value = 1.032295
$ws.Cells.Item(10, 9).Value = 0.09560560325825622
$ws.Cells.Item(10, 10).Value = 0.09560560325825622
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 7.407905
$ws.Cells.Item(10, 14).Value = 22.223715
$ws.Cells.Item(10, 15).Value = 0.1577242380174723
$ws.Cells.Item(10, 16).Value = 0.1577242380174723
$ws.Cells.Item(10, 17).Value = 2.549047763991667
$ws.Cells.Item(10, 18).Value = 22.941429875925
$ws.Cells.Item(10, 19).Value = 0.01507932092410923
$ws.Cells.Item(10, 20).Value = 0.01507932092410923

# Row 11
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "S100a8"
$ws.Cells.Item(11, 3).Value = "Tlr4"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 0.3333333333333333
$ws.Cells.Item(11, 7).Value = 0.3440983333333333
$ws.Cells.Item(11, 8).Value = 1.032295
$ws.Cells.Item(11, 9).Value = 0.09560560325825622
$ws.Cells.Item(11, 10).Value = 0.09560560325825622
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 16.48752133333333
$ws.Cells.Item(11, 14).Value = 49.462564
$ws.Cells.Item(11, 15).Value = 0.3510414535684271
$ws.Cells.Item(11, 16).Value = 0.3510414535684271
$ws.Cells.Item(11, 17).Value = 5.673328611597778
$ws.Cells.Item(11, 18).Value = 51.05995750438
$ws.Cells.Item(11, 19).Value = 0.03356152993706461
$ws.Cells.Item(11, 20).Value = 0.03356152993706462

# Row 12
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "S100a8"
$ws.Cells.Item(12, 3).Value = "Tlr4"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0.3333333333333333
$ws.Cells.Item(12, 7).Value = 0.3440983333333333
$ws.Cells.Item(12, 8).Value = 1.032295
$ws.Cells.Item(12, 9).Value = 0.09560560325825622
$ws.Cells.Item(12, 10).Value = 0.09560560325825622
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 5.464566666666666
$ws.Cells.Item(12, 14).Value = 16.3937
$ws.Cells.Item(12, 15).Value = 0.1163479571613943
$ws.Cells.Item(12, 16).Value = 0.1163479571613943
$ws.Cells.Item(12, 17).Value = 1.880348282388889
$ws.Cells.Item(12, 18).Value = 16.9231345415
$ws.Cells.Item(12, 19).Value = 0.01112351663228085
$ws.Cells.Item(12, 20).Value = 0.01112351663228085

# Row 13
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "S100a8"
$ws.Cells.Item(13, 3).Value = "Tlr4"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0.3333333333333333
$ws.Cells.Item(13, 7).Value = 0.3440983333333333
$ws.Cells.Item(13, 8).Value = 1.032295
$ws.Cells.Item(13, 9).Value = 0.09560560325825622
$ws.Cells.Item(13, 10).Value = 0.09560560325825622
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 17.60745533333333
$ws.Cells.Item(13, 14).Value = 52.822366
$ws.Cells.Item(13, 15).Value = 0.3748863512527063
$ws.Cells.Item(13, 16).Value = 0.3748863512527063
$ws.Cells.Item(13, 17).Value = 6.058696034441112
$ws.Cells.Item(13, 18).Value = 54.52826430997
$ws.Cells.Item(13, 19).Value = 0.03584123576480152
$ws.Cells.Item(13, 20).Value = 0.03584123576480153

# Row 14
$ws.Cells.Item(14, 1).Value = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value = "S100a8"
$ws.Cells.Item(14, 3).Value = "Tlr4"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 1.926427333333333
$ws.Cells.Item(14, 8).Value = 5.779282
$ws.Cells.Item(14, 9).Value = 0.5352459733017999
$ws.Cells.Item(14, 10).Value = 0.5352459733017999
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 7.407905
$ws.Cells.Item(14, 14).Value = 22.223715
$ws.Cells.Item(14, 15).Value = 0.1577242380174723
$ws.Cells.Item(14, 16).Value = 0.1577242380174723
$ws.Cells.Item(14, 17).Value = 14.27079067473667
$ws.Cells.Item(14, 18).Value = 128.43711607263
$ws.Cells.Item(14, 19).Value = 0.08442126329094671
$ws.Cells.Item(14, 20).Value = 0.08442126329094672

# Row 15
$ws.Cells.Item(15, 1).Value = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value = "S100a8"
$ws.Cells.Item(15, 3).Value = "Tlr4"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 1.926427333333333
$ws.Cells.Item(15, 8).Value = 5.779282
$ws.Cells.Item(15, 9).Value = 0.5352459733017999
$ws.Cells.Item(15, 10).Value = 0.5352459733017999
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 16.48752133333333
$ws.Cells.Item(15, 14).Value = 49.462564
$ws.Cells.Item(15, 15).Value = 0.3510414535684271
$ws.Cells.Item(15, 16).Value = 0.3510414535684271
$ws.Cells.Item(15, 17).Value = 31.76201175544978
$ws.Cells.Item(15, 18).Value = 285.858105799048
$ws.Cells.Item(15, 19).Value = 0.1878935244845114
$ws.Cells.Item(15, 20).Value = 0.1878935244845114

# Row 16
$ws.Cells.Item(16, 1).Value = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value = "S100a8"
$ws.Cells.Item(16, 3).Value = "Tlr4"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 1.926427333333333
$ws.Cells.Item(16, 8).Value = 5.779282
$ws.Cells.Item(16, 9).Value = 0.5352459733017999
$ws.Cells.Item(16, 10).Value = 0.5352459733017999
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 5.464566666666666
$ws.Cells.Item(16, 14).Value = 16.3937
$ws.Cells.Item(16, 15).Value = 0.1163479571613943
$ws.Cells.Item(16, 16).Value = 0.1163479571613943
$ws.Cells.Item(16, 17).Value = 10.52709059148889
$ws.Cells.Item(16, 18).Value = 94.7438153234
$ws.Cells.Item(16, 19).Value = 0.0622747755725266
$ws.Cells.Item(16, 20).Value = 0.06227477557252661

# Row 17
$ws.Cells.Item(17, 1).Value = "Resolving-Mac"
$ws.Cells.Item(17, 2).Value = "S100a8"
$ws.Cells.Item(17, 3).Value = "Tlr4"
$ws.Cells.Item(17, 4).Value = "Resolving-Mac"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 1.926427333333333
$ws.Cells.Item(17, 8).Value = 5.779282
$ws.Cells.Item(17, 9).Value = 0.5352459733017999
$ws.Cells.Item(17, 10).Value = 0.5352459733017999
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 17.60745533333333
$ws.Cells.Item(17, 14).Value = 52.822366
$ws.Cells.Item(17, 15).Value = 0.3748863512527063
$ws.Cells.Item(17, 16).Value = 0.3748863512527063
$ws.Cells.Item(17, 17).Value = 33.91948322457912
$ws.Cells.Item(17, 18).Value = 305.275349021212
$ws.Cells.Item(17, 19).Value = 0.2006564099538152
$ws.Cells.Item(17, 20).Value = 0.2006564099538152
